# Map a new "RunOptions" interface snippet into the Snippets table.
# This inserts a new row at position 283 (pushing the existing rows
# 283-410 down to 284-411) and fills it with the new snippet metadata,
# then grows the table/list object to cover the extra row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 283, shifting subsequent rows down.
$ws.Rows.Item(283).Insert()

# Populate the new row with the new snippet mapping.
$ws.Range("A283").Value = "Excel"
$ws.Range("B283").Value = "RunOptions"
$ws.Range("D283").Value = "interface"
$ws.Range("E283").Value = "excel-workbook-undo-grouping"
$ws.Range("F283").Value = "formatWithGrouping"

# Grow the "Snippets" table to include the newly inserted row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F411"))

# Reflect the author's final selection/scroll position in the sheet view.
[void]$ws.Range("F283").Select()
$excel.ActiveWindow.ScrollRow = 276
